# Item 12 had a typo ("REAME.md") and referenced a "push" step that no
# longer applies; replace the whole line with a pointer to the rest of
# the assignment on GitHub, and make the entire paragraph (text + the
# paragraph mark) bold.

$d = $word.ActiveDocument

$oldText = "12. push your first assignment to GitHub and add a REAME.md."
$newText = "12. The remainder of this assignment can be found at https://github.com/filiptosic/DSI_assignment2"

# Find the paragraph that still holds the old (three-run) text and keep a
# handle on the Paragraph object itself -- not just a Range -- so it stays
# valid once we rewrite the text inside it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq $oldText) {
        $target = $candidate
        break
    }
}

# Replace the runs' combined text with the new sentence, leaving the
# paragraph mark (end-of-paragraph `\r`) untouched.
$paraRange = $target.Range
$textOnly = $d.Range($paraRange.Start, $paraRange.End - 1)
$textOnly.Text = $newText

# Bold the new run and the paragraph mark itself, matching the source:
# <w:pPr><w:rPr><w:b/></w:rPr></w:pPr> + <w:r><w:rPr><w:b/>...
$target.Range.Font.Bold = 1
